# Update "想去人数" (interest count, column F) figures across sheets to
# reflect the latest scrape snapshot, as published to gh-pages.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 98
$ws1.Range("F6").Value  = 919
$ws1.Range("F8").Value  = 4869
$ws1.Range("F9").Value  = 4869
$ws1.Range("F10").Value = 113
$ws1.Range("F12").Value = 179
$ws1.Range("F14").Value = 206
$ws1.Range("F16").Value = 7863
$ws1.Range("F17").Value = 7863
$ws1.Range("F23").Value = 2272
$ws1.Range("F25").Value = 2501
$ws1.Range("F26").Value = 13
$ws1.Range("F28").Value = 6270
$ws1.Range("F33").Value = 455
$ws1.Range("F34").Value = 6639
$ws1.Range("F42").Value = 27
$ws1.Range("F43").Value = 2487
$ws1.Range("F50").Value = 57

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 242

# --- Sheet "全部类型" (All types, aggregate of the other sheets) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value  = 242
$ws4.Range("F7").Value  = 98
$ws4.Range("F10").Value = 4869
$ws4.Range("F11").Value = 4869
$ws4.Range("F12").Value = 113
$ws4.Range("F14").Value = 179
$ws4.Range("F17").Value = 7864
$ws4.Range("F18").Value = 7864
$ws4.Range("F26").Value = 2272
$ws4.Range("F27").Value = 2501
$ws4.Range("F28").Value = 13
$ws4.Range("F31").Value = 6270
$ws4.Range("F36").Value = 455
$ws4.Range("F37").Value = 6639
$ws4.Range("F42").Value = 27
$ws4.Range("F44").Value = 2487
$ws4.Range("F51").Value = 57
